# Updates the cryptos list ("cryptos.xlsx") Price and Volume(1h) columns
# (plus B/C for two coin rows that swapped ranking position) to match the
# latest scrape, matching commit:
#   "Updated cryptos list on Wed Sep 11 20:27:39 UTC 2024 with GitHub Actions"
#
# Price/Volume cells are stored as text (they use thousands separators like
# "57.591.04" and padded percentages like "  +0.21%  "), so every write uses
# a leading apostrophe to force text entry, then resets the cell style back
# to "Normal" so Excel does not silently apply a Text number format to
# values that would otherwise parse as numbers (e.g. "531.60" -> 531.6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.591.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.21%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.318.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.03%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'531.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.92%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'132.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.57%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.15%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.66%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.341.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.02%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.07%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.05%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'5.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.36%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.17%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.737.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.80%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'23.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.54%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'57.425.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.12%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.30%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.328.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.47%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'338.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.04%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -2.06%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.28%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -1.76%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.18%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'62.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.39%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.19%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'8.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.94%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.993"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.30%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.66%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'173.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.95%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +1.04%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0₃0726"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.01%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.04%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'18.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.28%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -4.29%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.917"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.92%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.58%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'OKB"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'39.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.40%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'Stacks"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'1.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.57%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +10.57%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'149.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.24%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -3.26%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'3.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.84%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'281.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.11%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0930"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.47%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0502"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.84%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'18.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.17%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.65%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.95%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +5.00%  "
$ws.Range("E51").Style = "Normal"
